# Added new categorical values for preparation medium and storage medium
# Closes #9
#
# - preparation_medium sheet: add Trumps fixative, DMEM, Biops buffer and
#   2% PFA/2.5% Glutaraldehyde; rename "Neutral Buffered Formalin (NBF)" to
#   "NBF (Neutral Buffered Formalin)".
# - storage_medium sheet: add Cyro-EM and 2% PFA/2.5% Glutaraldehyde; rename
#   "Neutral Buffered Formalin (NBF)" to "NBF (Neutral Buffered Formalin)".
# - bump the related dataValidation list ranges on the Sample Section sheet.
# - refresh the .metadata pav:createdOn timestamp.

$wb = $excel.ActiveWorkbook

$prepMedium = @(
    ,('NBF (Neutral Buffered Formalin)', 'http://purl.obolibrary.org/obo/OBIB_0000213')
    ,('Allprotect tissue reagent (ALL)', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000118')
    ,('CLARITY hydrogel', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000134')
    ,('Trumps fixative', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000331')
    ,('Inflated (OCT)', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000123')
    ,('DMEM', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185409')
    ,('PFA (Paraformaldehyde)', 'http://purl.obolibrary.org/obo/CHEBI_61538')
    ,('Fixed frozen OCT (Formalin, sucrose protected)', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000116')
    ,('Unknown', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998')
    ,('Fresh frozen OCT', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000126')
    ,('2% PFA/2.5% Glutaraldehyde', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000332')
    ,('Bouin''s', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000140')
    ,('Methanol', 'http://purl.obolibrary.org/obo/CHEBI_17790')
    ,('PAXgene tissue kit (PXT)', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185113')
    ,('PBS', 'http://purl.obolibrary.org/obo/OBI_0100046')
    ,('Ethanol', 'http://purl.obolibrary.org/obo/CHEBI_16236')
    ,('Inflated (Agarose)', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000106')
    ,('PLP (Periodate-Lysine-Paraformaldehyde)', 'http://purl.bioontology.org/ontology/MESH/C046311')
    ,('MACS tissue storage solution', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000105')
    ,('Fresh frozen CMC', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000130')
    ,('Fresh frozen gelatin', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000198')
    ,('RNAlater', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348')
    ,('Biops buffer', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000330')
    ,('Fixed frozen OCT (Cytofix/Cytoperm)', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000149')
    ,('None', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C41132')
    ,('Fixed frozen OCT (PFA, sucrose protected)', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000147')
)

$storageMedium = @(
    ,('PBS', 'http://purl.obolibrary.org/obo/OBI_0100046')
    ,('OCT', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63523')
    ,('NBF (Neutral Buffered Formalin)', 'http://purl.obolibrary.org/obo/OBIB_0000213')
    ,('Ethanol', 'http://purl.obolibrary.org/obo/CHEBI_16236')
    ,('Allprotect tissue reagent (ALL)', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000118')
    ,('DMSO (no serum)', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000115')
    ,('MACS tissue storage solution', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000105')
    ,('PFA (Paraformaldehyde)', 'http://purl.obolibrary.org/obo/CHEBI_61538')
    ,('Tris-EDTA', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000135')
    ,('Unknown', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998')
    ,('Gelatin', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C65802')
    ,('DMSO (serum)', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125')
    ,('RNAlater', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348')
    ,('Cyro-EM', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333')
    ,('FFPE (Paraffin embedded)', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C143028')
    ,('CMC', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C83594')
    ,('None', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C41132')
    ,('2% PFA/2.5% Glutaraldehyde', 'https://purl.humanatlas.io/vocab/hravs#HRAVS_0000332')
    ,('Methanol', 'http://purl.obolibrary.org/obo/CHEBI_17790')
    ,('PAXgene tissue kit (PXT)', 'http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185113')
)

# --- preparation_medium sheet: rewrite full A:B list (22 -> 26 rows) ---
$wsPrep = $wb.Worksheets.Item("preparation_medium")
for ($i = 0; $i -lt $prepMedium.Length; $i++) {
    $row = $i + 1
    $pair = $prepMedium[$i]
    $wsPrep.Cells.Item($row, 1).Value = $pair[0]
    $wsPrep.Cells.Item($row, 2).Value = $pair[1]
}

# --- storage_medium sheet: rewrite full A:B list (18 -> 20 rows) ---
$wsStorage = $wb.Worksheets.Item("storage_medium")
for ($i = 0; $i -lt $storageMedium.Length; $i++) {
    $row = $i + 1
    $pair = $storageMedium[$i]
    $wsStorage.Cells.Item($row, 1).Value = $pair[0]
    $wsStorage.Cells.Item($row, 2).Value = $pair[1]
}

# --- Sample Section sheet: bump the dataValidation list ranges that point
#     at the two lookup sheets whose row counts changed ---
$wsMain = $wb.Worksheets.Item("Sample Section")
$wsMain.Range("G2:G1001").Validation.Modify(3, 1, 1, "'preparation_medium'!`$A`$1:`$A`$26")
$wsMain.Range("K2:K1001").Validation.Modify(3, 1, 1, "'storage_medium'!`$A`$1:`$A`$20")

# --- .metadata sheet: refresh pav:createdOn timestamp ---
$wsMeta = $wb.Worksheets.Item(".metadata")
$wsMeta.Cells.Item(2, 3).Value = "2024-03-12T09:42:42-07:00"
